$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update the "Förändrad" (column C) date for all data rows 2-105 from 46072 to 46073 ---
for ($r = 2; $r -le 105; $r++) {
    $ws.Cells.Item($r, 3).Value = 46073
}

# --- Step 2: Apply the row-content permutation (rows were re-ordered in the source refresh) ---

# Row 6 -> becomes "A 60261-2024"
$ws.Range("A6").Value = "A 60261-2024"
$ws.Range("B6").Value = 45642
$ws.Range("G6").Value = 0.7
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Blåsippa"
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/artfynd/A 60261-2024 artfynd.xlsx", "A 60261-2024")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/kartor/A 60261-2024 karta.png", "A 60261-2024")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomål/A 60261-2024 FSC-klagomål.docx", "A 60261-2024")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomålsmail/A 60261-2024 FSC-klagomål mail.docx", "A 60261-2024")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsyn/A 60261-2024 tillsynsbegäran.docx", "A 60261-2024")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsynsmail/A 60261-2024 tillsynsbegäran mail.docx", "A 60261-2024")'

# Row 8 -> becomes "A 25789-2021"
$ws.Range("A8").Value = "A 25789-2021"
$ws.Range("B8").Value = 44343.8302662037
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Gullklöver"
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/artfynd/A 25789-2021 artfynd.xlsx", "A 25789-2021")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/kartor/A 25789-2021 karta.png", "A 25789-2021")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomål/A 25789-2021 FSC-klagomål.docx", "A 25789-2021")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomålsmail/A 25789-2021 FSC-klagomål mail.docx", "A 25789-2021")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsyn/A 25789-2021 tillsynsbegäran.docx", "A 25789-2021")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsynsmail/A 25789-2021 tillsynsbegäran mail.docx", "A 25789-2021")'

# Row 9 -> becomes "A 15108-2023"
$ws.Range("A9").Value = "A 15108-2023"
$ws.Range("B9").Value = 45016
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Skogsalm"
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/artfynd/A 15108-2023 artfynd.xlsx", "A 15108-2023")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/kartor/A 15108-2023 karta.png", "A 15108-2023")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomål/A 15108-2023 FSC-klagomål.docx", "A 15108-2023")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/klagomålsmail/A 15108-2023 FSC-klagomål mail.docx", "A 15108-2023")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsyn/A 15108-2023 tillsynsbegäran.docx", "A 15108-2023")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1484/tillsynsmail/A 15108-2023 tillsynsbegäran mail.docx", "A 15108-2023")'

# Row 24 -> becomes "A 6509-2025"
$ws.Range("A24").Value = "A 6509-2025"
$ws.Range("B24").Value = 45699.65163194444
$ws.Range("G24").Value = 4.3
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = ""
$ws.Range("S24").ClearContents()
$ws.Range("T24").ClearContents()
$ws.Range("V24").ClearContents()
$ws.Range("W24").ClearContents()
$ws.Range("X24").ClearContents()
$ws.Range("Y24").ClearContents()

# Row 25 -> becomes "A 60262-2024"
$ws.Range("A25").Value = "A 60262-2024"
$ws.Range("B25").Value = 45642
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = ""
$ws.Range("S25").ClearContents()
$ws.Range("T25").ClearContents()
$ws.Range("V25").ClearContents()
$ws.Range("W25").ClearContents()
$ws.Range("X25").ClearContents()
$ws.Range("Y25").ClearContents()

# Row 26 -> becomes "A 60265-2024"
$ws.Range("A26").Value = "A 60265-2024"
$ws.Range("B26").Value = 45642
$ws.Range("G26").Value = 1.7
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = ""
$ws.Range("S26").ClearContents()
$ws.Range("T26").ClearContents()
$ws.Range("V26").ClearContents()
$ws.Range("W26").ClearContents()
$ws.Range("X26").ClearContents()
$ws.Range("Y26").ClearContents()

# Row 27 -> becomes "A 53423-2023"
$ws.Range("A27").Value = "A 53423-2023"
$ws.Range("B27").Value = 45229
$ws.Range("G27").Value = 1.5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ""
$ws.Range("S27").ClearContents()
$ws.Range("T27").ClearContents()
$ws.Range("V27").ClearContents()
$ws.Range("W27").ClearContents()
$ws.Range("X27").ClearContents()
$ws.Range("Y27").ClearContents()

# Row 28 -> becomes "A 48148-2025"
$ws.Range("A28").Value = "A 48148-2025"
$ws.Range("B28").Value = 45933.40649305555
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ""
$ws.Range("S28").ClearContents()
$ws.Range("T28").ClearContents()
$ws.Range("V28").ClearContents()
$ws.Range("W28").ClearContents()
$ws.Range("X28").ClearContents()
$ws.Range("Y28").ClearContents()

# Row 29 -> becomes "A 48153-2025"
$ws.Range("A29").Value = "A 48153-2025"
$ws.Range("B29").Value = 45933.41109953704
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = ""
$ws.Range("S29").ClearContents()
$ws.Range("T29").ClearContents()
$ws.Range("V29").ClearContents()
$ws.Range("W29").ClearContents()
$ws.Range("X29").ClearContents()
$ws.Range("Y29").ClearContents()

# Row 30 -> becomes "A 60255-2024"
$ws.Range("A30").Value = "A 60255-2024"
$ws.Range("B30").Value = 45642
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = ""
$ws.Range("S30").ClearContents()
$ws.Range("T30").ClearContents()
$ws.Range("V30").ClearContents()
$ws.Range("W30").ClearContents()
$ws.Range("X30").ClearContents()
$ws.Range("Y30").ClearContents()

# Row 31 -> becomes "A 52-2024"
$ws.Range("A31").Value = "A 52-2024"
$ws.Range("B31").Value = 45293
$ws.Range("G31").Value = 2.7
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = ""
$ws.Range("S31").ClearContents()
$ws.Range("T31").ClearContents()
$ws.Range("V31").ClearContents()
$ws.Range("W31").ClearContents()
$ws.Range("X31").ClearContents()
$ws.Range("Y31").ClearContents()

# Row 32 -> becomes "A 10579-2024"
$ws.Range("A32").Value = "A 10579-2024"
$ws.Range("B32").Value = 45366.66769675926
$ws.Range("G32").Value = 1.3
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = ""
$ws.Range("S32").ClearContents()
$ws.Range("T32").ClearContents()
$ws.Range("V32").ClearContents()
$ws.Range("W32").ClearContents()
$ws.Range("X32").ClearContents()
$ws.Range("Y32").ClearContents()

# Row 33 -> becomes "A 40154-2025"
$ws.Range("A33").Value = "A 40154-2025"
$ws.Range("B33").Value = 45894.5980787037
$ws.Range("G33").Value = 1.7
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = ""
$ws.Range("S33").ClearContents()
$ws.Range("T33").ClearContents()
$ws.Range("V33").ClearContents()
$ws.Range("W33").ClearContents()
$ws.Range("X33").ClearContents()
$ws.Range("Y33").ClearContents()

# Row 34 -> becomes "A 22658-2025"
$ws.Range("A34").Value = "A 22658-2025"
$ws.Range("B34").Value = 45789.47318287037
$ws.Range("G34").Value = 4.5
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = ""
$ws.Range("S34").ClearContents()
$ws.Range("T34").ClearContents()
$ws.Range("V34").ClearContents()
$ws.Range("W34").ClearContents()
$ws.Range("X34").ClearContents()
$ws.Range("Y34").ClearContents()

# Row 35 -> becomes "A 40683-2025"
$ws.Range("A35").Value = "A 40683-2025"
$ws.Range("B35").Value = 45896
$ws.Range("G35").Value = 0.9
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = ""
$ws.Range("S35").ClearContents()
$ws.Range("T35").ClearContents()
$ws.Range("V35").ClearContents()
$ws.Range("W35").ClearContents()
$ws.Range("X35").ClearContents()
$ws.Range("Y35").ClearContents()

# Row 36 -> becomes "A 40348-2025"
$ws.Range("A36").Value = "A 40348-2025"
$ws.Range("B36").Value = 45895.480625
$ws.Range("G36").Value = 1.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = ""
$ws.Range("S36").ClearContents()
$ws.Range("T36").ClearContents()
$ws.Range("V36").ClearContents()
$ws.Range("W36").ClearContents()
$ws.Range("X36").ClearContents()
$ws.Range("Y36").ClearContents()

# Row 37 -> becomes "A 13246-2024"
$ws.Range("A37").Value = "A 13246-2024"
$ws.Range("B37").Value = 45386.56986111111
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = ""
$ws.Range("S37").ClearContents()
$ws.Range("T37").ClearContents()
$ws.Range("V37").ClearContents()
$ws.Range("W37").ClearContents()
$ws.Range("X37").ClearContents()
$ws.Range("Y37").ClearContents()

# Row 38 -> becomes "A 41381-2025"
$ws.Range("A38").Value = "A 41381-2025"
$ws.Range("B38").Value = 45898
$ws.Range("G38").Value = 1.2
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = ""
$ws.Range("S38").ClearContents()
$ws.Range("T38").ClearContents()
$ws.Range("V38").ClearContents()
$ws.Range("W38").ClearContents()
$ws.Range("X38").ClearContents()
$ws.Range("Y38").ClearContents()

# Row 39 -> becomes "A 49105-2025"
$ws.Range("A39").Value = "A 49105-2025"
$ws.Range("B39").Value = 45937.72574074074
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = ""
$ws.Range("S39").ClearContents()
$ws.Range("T39").ClearContents()
$ws.Range("V39").ClearContents()
$ws.Range("W39").ClearContents()
$ws.Range("X39").ClearContents()
$ws.Range("Y39").ClearContents()

# Row 40 -> becomes "A 22605-2025"
$ws.Range("A40").Value = "A 22605-2025"
$ws.Range("B40").Value = 45789
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = ""
$ws.Range("S40").ClearContents()
$ws.Range("T40").ClearContents()
$ws.Range("V40").ClearContents()
$ws.Range("W40").ClearContents()
$ws.Range("X40").ClearContents()
$ws.Range("Y40").ClearContents()

# Row 41 -> becomes "A 41395-2025"
$ws.Range("A41").Value = "A 41395-2025"
$ws.Range("B41").Value = 45901.31965277778
$ws.Range("G41").Value = 1.9
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = ""
$ws.Range("S41").ClearContents()
$ws.Range("T41").ClearContents()
$ws.Range("V41").ClearContents()
$ws.Range("W41").ClearContents()
$ws.Range("X41").ClearContents()
$ws.Range("Y41").ClearContents()

# Row 42 -> becomes "A 34623-2025"
$ws.Range("A42").Value = "A 34623-2025"
$ws.Range("B42").Value = 45847
$ws.Range("G42").Value = 3.4
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = ""
$ws.Range("S42").ClearContents()
$ws.Range("T42").ClearContents()
$ws.Range("V42").ClearContents()
$ws.Range("W42").ClearContents()
$ws.Range("X42").ClearContents()
$ws.Range("Y42").ClearContents()

# Row 43 -> becomes "A 58619-2024"
$ws.Range("A43").Value = "A 58619-2024"
$ws.Range("B43").Value = 45635.56216435185
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = ""
$ws.Range("S43").ClearContents()
$ws.Range("T43").ClearContents()
$ws.Range("V43").ClearContents()
$ws.Range("W43").ClearContents()
$ws.Range("X43").ClearContents()
$ws.Range("Y43").ClearContents()

# Row 45 -> becomes "A 41581-2023"
$ws.Range("A45").Value = "A 41581-2023"
$ws.Range("B45").Value = 45173
$ws.Range("G45").Value = 2.3
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = ""
$ws.Range("S45").ClearContents()
$ws.Range("T45").ClearContents()
$ws.Range("V45").ClearContents()
$ws.Range("W45").ClearContents()
$ws.Range("X45").ClearContents()
$ws.Range("Y45").ClearContents()

# Row 46 -> becomes "A 19794-2023"
$ws.Range("A46").Value = "A 19794-2023"
$ws.Range("B46").Value = 45051
$ws.Range("G46").Value = 2.2
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = ""
$ws.Range("S46").ClearContents()
$ws.Range("T46").ClearContents()
$ws.Range("V46").ClearContents()
$ws.Range("W46").ClearContents()
$ws.Range("X46").ClearContents()
$ws.Range("Y46").ClearContents()

# Row 47 -> becomes "A 63886-2023"
$ws.Range("A47").Value = "A 63886-2023"
$ws.Range("B47").Value = 45278
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = ""
$ws.Range("S47").ClearContents()
$ws.Range("T47").ClearContents()
$ws.Range("V47").ClearContents()
$ws.Range("W47").ClearContents()
$ws.Range("X47").ClearContents()
$ws.Range("Y47").ClearContents()

# Row 48 -> becomes "A 53191-2025"
$ws.Range("A48").Value = "A 53191-2025"
$ws.Range("B48").Value = 45958.56690972222
$ws.Range("G48").Value = 0.9
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = ""
$ws.Range("S48").ClearContents()
$ws.Range("T48").ClearContents()
$ws.Range("V48").ClearContents()
$ws.Range("W48").ClearContents()
$ws.Range("X48").ClearContents()
$ws.Range("Y48").ClearContents()

# Row 49 -> becomes "A 43151-2025"
$ws.Range("A49").Value = "A 43151-2025"
$ws.Range("B49").Value = 45910.33892361111
$ws.Range("G49").Value = 2.9
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = ""
$ws.Range("S49").ClearContents()
$ws.Range("T49").ClearContents()
$ws.Range("V49").ClearContents()
$ws.Range("W49").ClearContents()
$ws.Range("X49").ClearContents()
$ws.Range("Y49").ClearContents()

# Row 50 -> becomes "A 25395-2025"
$ws.Range("A50").Value = "A 25395-2025"
$ws.Range("B50").Value = 45800.64942129629
$ws.Range("G50").Value = 1.7
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = ""
$ws.Range("S50").ClearContents()
$ws.Range("T50").ClearContents()
$ws.Range("V50").ClearContents()
$ws.Range("W50").ClearContents()
$ws.Range("X50").ClearContents()
$ws.Range("Y50").ClearContents()

# Row 51 -> becomes "A 44022-2025"
$ws.Range("A51").Value = "A 44022-2025"
$ws.Range("B51").Value = 45915.47851851852
$ws.Range("G51").Value = 2.1
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = ""
$ws.Range("S51").ClearContents()
$ws.Range("T51").ClearContents()
$ws.Range("V51").ClearContents()
$ws.Range("W51").ClearContents()
$ws.Range("X51").ClearContents()
$ws.Range("Y51").ClearContents()

# Row 52 -> becomes "A 54366-2025"
$ws.Range("A52").Value = "A 54366-2025"
$ws.Range("B52").Value = 45965.40274305556
$ws.Range("G52").Value = 2.6
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ""
$ws.Range("S52").ClearContents()
$ws.Range("T52").ClearContents()
$ws.Range("V52").ClearContents()
$ws.Range("W52").ClearContents()
$ws.Range("X52").ClearContents()
$ws.Range("Y52").ClearContents()

# Row 53 -> becomes "A 48403-2024"
$ws.Range("A53").Value = "A 48403-2024"
$ws.Range("B53").Value = 45590.61787037037
$ws.Range("G53").Value = 1.1
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = ""
$ws.Range("S53").ClearContents()
$ws.Range("T53").ClearContents()
$ws.Range("V53").ClearContents()
$ws.Range("W53").ClearContents()
$ws.Range("X53").ClearContents()
$ws.Range("Y53").ClearContents()

# Row 54 -> becomes "A 20685-2024"
$ws.Range("A54").Value = "A 20685-2024"
$ws.Range("B54").Value = 45436.58555555555
$ws.Range("G54").Value = 1.5
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = ""
$ws.Range("S54").ClearContents()
$ws.Range("T54").ClearContents()
$ws.Range("V54").ClearContents()
$ws.Range("W54").ClearContents()
$ws.Range("X54").ClearContents()
$ws.Range("Y54").ClearContents()

# Row 55 -> becomes "A 37149-2023"
$ws.Range("A55").Value = "A 37149-2023"
$ws.Range("B55").Value = 45155
$ws.Range("G55").Value = 1.5
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = ""
$ws.Range("S55").ClearContents()
$ws.Range("T55").ClearContents()
$ws.Range("V55").ClearContents()
$ws.Range("W55").ClearContents()
$ws.Range("X55").ClearContents()
$ws.Range("Y55").ClearContents()

# Row 56 -> becomes "A 57532-2025"
$ws.Range("A56").Value = "A 57532-2025"
$ws.Range("B56").Value = 45980.63876157408
$ws.Range("G56").Value = 2.5
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = ""
$ws.Range("S56").ClearContents()
$ws.Range("T56").ClearContents()
$ws.Range("V56").ClearContents()
$ws.Range("W56").ClearContents()
$ws.Range("X56").ClearContents()
$ws.Range("Y56").ClearContents()

# Row 57 -> becomes "A 57530-2025"
$ws.Range("A57").Value = "A 57530-2025"
$ws.Range("B57").Value = 45980.63579861111
$ws.Range("G57").Value = 2.9
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("O57").Value = 0
$ws.Range("P57").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("R57").Value = ""
$ws.Range("S57").ClearContents()
$ws.Range("T57").ClearContents()
$ws.Range("V57").ClearContents()
$ws.Range("W57").ClearContents()
$ws.Range("X57").ClearContents()
$ws.Range("Y57").ClearContents()

# Row 58 -> becomes "A 7149-2025"
$ws.Range("A58").Value = "A 7149-2025"
$ws.Range("B58").Value = 45702.37914351852
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("N58").Value = 0
$ws.Range("O58").Value = 0
$ws.Range("P58").Value = 0
$ws.Range("Q58").Value = 0
$ws.Range("R58").Value = ""
$ws.Range("S58").ClearContents()
$ws.Range("T58").ClearContents()
$ws.Range("V58").ClearContents()
$ws.Range("W58").ClearContents()
$ws.Range("X58").ClearContents()
$ws.Range("Y58").ClearContents()

# Row 59 -> becomes "A 43028-2023"
$ws.Range("A59").Value = "A 43028-2023"
$ws.Range("B59").Value = 45182
$ws.Range("G59").Value = 0.6
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("O59").Value = 0
$ws.Range("P59").Value = 0
$ws.Range("Q59").Value = 0
$ws.Range("R59").Value = ""
$ws.Range("S59").ClearContents()
$ws.Range("T59").ClearContents()
$ws.Range("V59").ClearContents()
$ws.Range("W59").ClearContents()
$ws.Range("X59").ClearContents()
$ws.Range("Y59").ClearContents()

# Row 61 -> becomes "A 35139-2025"
$ws.Range("A61").Value = "A 35139-2025"
$ws.Range("B61").Value = 45852.64978009259
$ws.Range("G61").Value = 0.7
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("O61").Value = 0
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("R61").Value = ""
$ws.Range("S61").ClearContents()
$ws.Range("T61").ClearContents()
$ws.Range("V61").ClearContents()
$ws.Range("W61").ClearContents()
$ws.Range("X61").ClearContents()
$ws.Range("Y61").ClearContents()

# Row 62 -> becomes "A 1265-2026"
$ws.Range("A62").Value = "A 1265-2026"
$ws.Range("B62").Value = 46031.46165509259
$ws.Range("G62").Value = 1.4
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("O62").Value = 0
$ws.Range("P62").Value = 0
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = ""
$ws.Range("S62").ClearContents()
$ws.Range("T62").ClearContents()
$ws.Range("V62").ClearContents()
$ws.Range("W62").ClearContents()
$ws.Range("X62").ClearContents()
$ws.Range("Y62").ClearContents()

# Row 63 -> becomes "A 35135-2025"
$ws.Range("A63").Value = "A 35135-2025"
$ws.Range("B63").Value = 45852.64123842592
$ws.Range("G63").Value = 0.1
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("O63").Value = 0
$ws.Range("P63").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("R63").Value = ""
$ws.Range("S63").ClearContents()
$ws.Range("T63").ClearContents()
$ws.Range("V63").ClearContents()
$ws.Range("W63").ClearContents()
$ws.Range("X63").ClearContents()
$ws.Range("Y63").ClearContents()

# Row 64 -> becomes "A 35640-2025"
$ws.Range("A64").Value = "A 35640-2025"
$ws.Range("B64").Value = 45859.47606481481
$ws.Range("G64").Value = 1.8
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("O64").Value = 0
$ws.Range("P64").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = ""
$ws.Range("S64").ClearContents()
$ws.Range("T64").ClearContents()
$ws.Range("V64").ClearContents()
$ws.Range("W64").ClearContents()
$ws.Range("X64").ClearContents()
$ws.Range("Y64").ClearContents()

# Row 65 -> becomes "A 2122-2025"
$ws.Range("A65").Value = "A 2122-2025"
$ws.Range("B65").Value = 45672.64579861111
$ws.Range("G65").Value = 1.2
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("O65").Value = 0
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 0
$ws.Range("R65").Value = ""
$ws.Range("S65").ClearContents()
$ws.Range("T65").ClearContents()
$ws.Range("V65").ClearContents()
$ws.Range("W65").ClearContents()
$ws.Range("X65").ClearContents()
$ws.Range("Y65").ClearContents()

# Row 66 -> becomes "A 8469-2023"
$ws.Range("A66").Value = "A 8469-2023"
$ws.Range("B66").Value = 44977.50028935185
$ws.Range("G66").Value = 1.1
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = ""
$ws.Range("S66").ClearContents()
$ws.Range("T66").ClearContents()
$ws.Range("V66").ClearContents()
$ws.Range("W66").ClearContents()
$ws.Range("X66").ClearContents()
$ws.Range("Y66").ClearContents()

# Row 67 -> becomes "A 38036-2025"
$ws.Range("A67").Value = "A 38036-2025"
$ws.Range("B67").Value = 45881
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("P67").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = ""
$ws.Range("S67").ClearContents()
$ws.Range("T67").ClearContents()
$ws.Range("V67").ClearContents()
$ws.Range("W67").ClearContents()
$ws.Range("X67").ClearContents()
$ws.Range("Y67").ClearContents()

# Row 68 -> becomes "A 2058-2025"
$ws.Range("A68").Value = "A 2058-2025"
$ws.Range("B68").Value = 45671
$ws.Range("G68").Value = 3.4
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("P68").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("R68").Value = ""
$ws.Range("S68").ClearContents()
$ws.Range("T68").ClearContents()
$ws.Range("V68").ClearContents()
$ws.Range("W68").ClearContents()
$ws.Range("X68").ClearContents()
$ws.Range("Y68").ClearContents()

# Row 69 -> becomes "A 25657-2021"
$ws.Range("A69").Value = "A 25657-2021"
$ws.Range("B69").Value = 44343.55598379629
$ws.Range("G69").Value = 9.800000000000001
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("O69").Value = 0
$ws.Range("P69").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("R69").Value = ""
$ws.Range("S69").ClearContents()
$ws.Range("T69").ClearContents()
$ws.Range("V69").ClearContents()
$ws.Range("W69").ClearContents()
$ws.Range("X69").ClearContents()
$ws.Range("Y69").ClearContents()

# Row 70 -> becomes "A 61985-2024"
$ws.Range("A70").Value = "A 61985-2024"
$ws.Range("B70").Value = 45653.67324074074
$ws.Range("G70").Value = 1.8
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("O70").Value = 0
$ws.Range("P70").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("R70").Value = ""
$ws.Range("S70").ClearContents()
$ws.Range("T70").ClearContents()
$ws.Range("V70").ClearContents()
$ws.Range("W70").ClearContents()
$ws.Range("X70").ClearContents()
$ws.Range("Y70").ClearContents()

# Row 71 -> becomes "A 3588-2026"
$ws.Range("A71").Value = "A 3588-2026"
$ws.Range("B71").Value = 46042.75251157407
$ws.Range("G71").Value = 1.5
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("O71").Value = 0
$ws.Range("P71").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("R71").Value = ""
$ws.Range("S71").ClearContents()
$ws.Range("T71").ClearContents()
$ws.Range("V71").ClearContents()
$ws.Range("W71").ClearContents()
$ws.Range("X71").ClearContents()
$ws.Range("Y71").ClearContents()

# Row 72 -> becomes "A 7509-2026"
$ws.Range("A72").Value = "A 7509-2026"
$ws.Range("B72").Value = 46059.58040509259
$ws.Range("G72").Value = 0.6
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("O72").Value = 0
$ws.Range("P72").Value = 0
$ws.Range("Q72").Value = 0
$ws.Range("R72").Value = ""
$ws.Range("S72").ClearContents()
$ws.Range("T72").ClearContents()
$ws.Range("V72").ClearContents()
$ws.Range("W72").ClearContents()
$ws.Range("X72").ClearContents()
$ws.Range("Y72").ClearContents()

# Row 73 -> becomes "A 7506-2026"
$ws.Range("A73").Value = "A 7506-2026"
$ws.Range("B73").Value = 46059.57417824074
$ws.Range("G73").Value = 2.4
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("O73").Value = 0
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = ""
$ws.Range("S73").ClearContents()
$ws.Range("T73").ClearContents()
$ws.Range("V73").ClearContents()
$ws.Range("W73").ClearContents()
$ws.Range("X73").ClearContents()
$ws.Range("Y73").ClearContents()

# Row 74 -> becomes "A 58793-2022"
$ws.Range("A74").Value = "A 58793-2022"
$ws.Range("B74").Value = 44903.42074074074
$ws.Range("G74").Value = 1.5
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("O74").Value = 0
$ws.Range("P74").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = ""
$ws.Range("S74").ClearContents()
$ws.Range("T74").ClearContents()
$ws.Range("V74").ClearContents()
$ws.Range("W74").ClearContents()
$ws.Range("X74").ClearContents()
$ws.Range("Y74").ClearContents()

# Row 75 -> becomes "A 8328-2026"
$ws.Range("A75").Value = "A 8328-2026"
$ws.Range("B75").Value = 46064.62018518519
$ws.Range("G75").Value = 1.3
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 0
$ws.Range("P75").Value = 0
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = ""
$ws.Range("S75").ClearContents()
$ws.Range("T75").ClearContents()
$ws.Range("V75").ClearContents()
$ws.Range("W75").ClearContents()
$ws.Range("X75").ClearContents()
$ws.Range("Y75").ClearContents()

# Row 76 -> becomes "A 9149-2026"
$ws.Range("A76").Value = "A 9149-2026"
$ws.Range("B76").Value = 46070.31185185185
$ws.Range("G76").Value = 3.4
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 0
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = ""
$ws.Range("S76").ClearContents()
$ws.Range("T76").ClearContents()
$ws.Range("V76").ClearContents()
$ws.Range("W76").ClearContents()
$ws.Range("X76").ClearContents()
$ws.Range("Y76").ClearContents()

# Row 77 -> becomes "A 25009-2023"
$ws.Range("A77").Value = "A 25009-2023"
$ws.Range("B77").Value = 45085
$ws.Range("G77").Value = 2.2
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("O77").Value = 0
$ws.Range("P77").Value = 0
$ws.Range("Q77").Value = 0
$ws.Range("R77").Value = ""
$ws.Range("S77").ClearContents()
$ws.Range("T77").ClearContents()
$ws.Range("V77").ClearContents()
$ws.Range("W77").ClearContents()
$ws.Range("X77").ClearContents()
$ws.Range("Y77").ClearContents()

# Row 78 -> becomes "A 44536-2024"
$ws.Range("A78").Value = "A 44536-2024"
$ws.Range("B78").Value = 45574.45354166667
$ws.Range("G78").Value = 1.7
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0
$ws.Range("R78").Value = ""
$ws.Range("S78").ClearContents()
$ws.Range("T78").ClearContents()
$ws.Range("V78").ClearContents()
$ws.Range("W78").ClearContents()
$ws.Range("X78").ClearContents()
$ws.Range("Y78").ClearContents()

# Row 79 -> becomes "A 37583-2024"
$ws.Range("A79").Value = "A 37583-2024"
$ws.Range("B79").Value = 45541.45914351852
$ws.Range("G79").Value = 0.7
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("P79").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = ""
$ws.Range("S79").ClearContents()
$ws.Range("T79").ClearContents()
$ws.Range("V79").ClearContents()
$ws.Range("W79").ClearContents()
$ws.Range("X79").ClearContents()
$ws.Range("Y79").ClearContents()

# Row 80 -> becomes "A 2053-2025"
$ws.Range("A80").Value = "A 2053-2025"
$ws.Range("B80").Value = 45671
$ws.Range("G80").Value = 0.9
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("O80").Value = 0
$ws.Range("P80").Value = 0
$ws.Range("Q80").Value = 0
$ws.Range("R80").Value = ""
$ws.Range("S80").ClearContents()
$ws.Range("T80").ClearContents()
$ws.Range("V80").ClearContents()
$ws.Range("W80").ClearContents()
$ws.Range("X80").ClearContents()
$ws.Range("Y80").ClearContents()

# Row 81 -> becomes "A 61991-2024"
$ws.Range("A81").Value = "A 61991-2024"
$ws.Range("B81").Value = 45653.70577546296
$ws.Range("G81").Value = 0.6
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("O81").Value = 0
$ws.Range("P81").Value = 0
$ws.Range("Q81").Value = 0
$ws.Range("R81").Value = ""
$ws.Range("S81").ClearContents()
$ws.Range("T81").ClearContents()
$ws.Range("V81").ClearContents()
$ws.Range("W81").ClearContents()
$ws.Range("X81").ClearContents()
$ws.Range("Y81").ClearContents()

# Row 82 -> becomes "A 37050-2023"
$ws.Range("A82").Value = "A 37050-2023"
$ws.Range("B82").Value = 45155.42420138889
$ws.Range("G82").Value = 4.4
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("O82").Value = 0
$ws.Range("P82").Value = 0
$ws.Range("Q82").Value = 0
$ws.Range("R82").Value = ""
$ws.Range("S82").ClearContents()
$ws.Range("T82").ClearContents()
$ws.Range("V82").ClearContents()
$ws.Range("W82").ClearContents()
$ws.Range("X82").ClearContents()
$ws.Range("Y82").ClearContents()

# Row 83 -> becomes "A 37077-2023"
$ws.Range("A83").Value = "A 37077-2023"
$ws.Range("B83").Value = 45155.48460648148
$ws.Range("G83").Value = 0.4
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("O83").Value = 0
$ws.Range("P83").Value = 0
$ws.Range("Q83").Value = 0
$ws.Range("R83").Value = ""
$ws.Range("S83").ClearContents()
$ws.Range("T83").ClearContents()
$ws.Range("V83").ClearContents()
$ws.Range("W83").ClearContents()
$ws.Range("X83").ClearContents()
$ws.Range("Y83").ClearContents()

# Row 84 -> becomes "A 44146-2023"
$ws.Range("A84").Value = "A 44146-2023"
$ws.Range("B84").Value = 45188
$ws.Range("G84").Value = 2.3
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("O84").Value = 0
$ws.Range("P84").Value = 0
$ws.Range("Q84").Value = 0
$ws.Range("R84").Value = ""
$ws.Range("S84").ClearContents()
$ws.Range("T84").ClearContents()
$ws.Range("V84").ClearContents()
$ws.Range("W84").ClearContents()
$ws.Range("X84").ClearContents()
$ws.Range("Y84").ClearContents()

# Row 85 -> becomes "A 58789-2022"
$ws.Range("A85").Value = "A 58789-2022"
$ws.Range("B85").Value = 44903.41709490741
$ws.Range("G85").Value = 0.4
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("O85").Value = 0
$ws.Range("P85").Value = 0
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = ""
$ws.Range("S85").ClearContents()
$ws.Range("T85").ClearContents()
$ws.Range("V85").ClearContents()
$ws.Range("W85").ClearContents()
$ws.Range("X85").ClearContents()
$ws.Range("Y85").ClearContents()

# Row 86 -> becomes "A 58794-2022"
$ws.Range("A86").Value = "A 58794-2022"
$ws.Range("B86").Value = 44903.42262731482
$ws.Range("G86").Value = 0.6
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("R86").Value = ""
$ws.Range("S86").ClearContents()
$ws.Range("T86").ClearContents()
$ws.Range("V86").ClearContents()
$ws.Range("W86").ClearContents()
$ws.Range("X86").ClearContents()
$ws.Range("Y86").ClearContents()

# Row 87 -> becomes "A 47304-2022"
$ws.Range("A87").Value = "A 47304-2022"
$ws.Range("B87").Value = 44853
$ws.Range("G87").Value = 2.8
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("O87").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = ""
$ws.Range("S87").ClearContents()
$ws.Range("T87").ClearContents()
$ws.Range("V87").ClearContents()
$ws.Range("W87").ClearContents()
$ws.Range("X87").ClearContents()
$ws.Range("Y87").ClearContents()

# Row 88 -> becomes "A 20683-2024"
$ws.Range("A88").Value = "A 20683-2024"
$ws.Range("B88").Value = 45436.58328703704
$ws.Range("G88").Value = 0.9
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("O88").Value = 0
$ws.Range("P88").Value = 0
$ws.Range("Q88").Value = 0
$ws.Range("R88").Value = ""
$ws.Range("S88").ClearContents()
$ws.Range("T88").ClearContents()
$ws.Range("V88").ClearContents()
$ws.Range("W88").ClearContents()
$ws.Range("X88").ClearContents()
$ws.Range("Y88").ClearContents()

# Row 89 -> becomes "A 59223-2022"
$ws.Range("A89").Value = "A 59223-2022"
$ws.Range("B89").Value = 44899
$ws.Range("G89").Value = 1.4
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("O89").Value = 0
$ws.Range("P89").Value = 0
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = ""
$ws.Range("S89").ClearContents()
$ws.Range("T89").ClearContents()
$ws.Range("V89").ClearContents()
$ws.Range("W89").ClearContents()
$ws.Range("X89").ClearContents()
$ws.Range("Y89").ClearContents()

# Row 90 -> becomes "A 8464-2023"
$ws.Range("A90").Value = "A 8464-2023"
$ws.Range("B90").Value = 44977.49489583333
$ws.Range("G90").Value = 0.5
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("O90").Value = 0
$ws.Range("P90").Value = 0
$ws.Range("Q90").Value = 0
$ws.Range("R90").Value = ""
$ws.Range("S90").ClearContents()
$ws.Range("T90").ClearContents()
$ws.Range("V90").ClearContents()
$ws.Range("W90").ClearContents()
$ws.Range("X90").ClearContents()
$ws.Range("Y90").ClearContents()

# Row 91 -> becomes "A 60969-2021"
$ws.Range("A91").Value = "A 60969-2021"
$ws.Range("B91").Value = 44497.5580787037
$ws.Range("G91").Value = 1.4
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("O91").Value = 0
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 0
$ws.Range("R91").Value = ""
$ws.Range("S91").ClearContents()
$ws.Range("T91").ClearContents()
$ws.Range("V91").ClearContents()
$ws.Range("W91").ClearContents()
$ws.Range("X91").ClearContents()
$ws.Range("Y91").ClearContents()

# Row 93 -> becomes "A 13239-2024"
$ws.Range("A93").Value = "A 13239-2024"
$ws.Range("B93").Value = 45386.5583449074
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("O93").Value = 0
$ws.Range("P93").Value = 0
$ws.Range("Q93").Value = 0
$ws.Range("R93").Value = ""
$ws.Range("S93").ClearContents()
$ws.Range("T93").ClearContents()
$ws.Range("V93").ClearContents()
$ws.Range("W93").ClearContents()
$ws.Range("X93").ClearContents()
$ws.Range("Y93").ClearContents()

# Row 94 -> becomes "A 39085-2023"
$ws.Range("A94").Value = "A 39085-2023"
$ws.Range("B94").Value = 45164.34875
$ws.Range("G94").Value = 0.6
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 0
$ws.Range("O94").Value = 0
$ws.Range("P94").Value = 0
$ws.Range("Q94").Value = 0
$ws.Range("R94").Value = ""
$ws.Range("S94").ClearContents()
$ws.Range("T94").ClearContents()
$ws.Range("V94").ClearContents()
$ws.Range("W94").ClearContents()
$ws.Range("X94").ClearContents()
$ws.Range("Y94").ClearContents()

# Row 95 -> becomes "A 440-2023"
$ws.Range("A95").Value = "A 440-2023"
$ws.Range("B95").Value = 44929.70428240741
$ws.Range("G95").Value = 4.2
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("O95").Value = 0
$ws.Range("P95").Value = 0
$ws.Range("Q95").Value = 0
$ws.Range("R95").Value = ""
$ws.Range("S95").ClearContents()
$ws.Range("T95").ClearContents()
$ws.Range("V95").ClearContents()
$ws.Range("W95").ClearContents()
$ws.Range("X95").ClearContents()
$ws.Range("Y95").ClearContents()

# Row 96 -> becomes "A 8474-2023"
$ws.Range("A96").Value = "A 8474-2023"
$ws.Range("B96").Value = 44977.50395833333
$ws.Range("G96").Value = 0.3
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = 0
$ws.Range("Q96").Value = 0
$ws.Range("R96").Value = ""
$ws.Range("S96").ClearContents()
$ws.Range("T96").ClearContents()
$ws.Range("V96").ClearContents()
$ws.Range("W96").ClearContents()
$ws.Range("X96").ClearContents()
$ws.Range("Y96").ClearContents()

# Row 97 -> becomes "A 18138-2025"
$ws.Range("A97").Value = "A 18138-2025"
$ws.Range("B97").Value = 45761.60342592592
$ws.Range("G97").Value = 0.8
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 0
$ws.Range("O97").Value = 0
$ws.Range("P97").Value = 0
$ws.Range("Q97").Value = 0
$ws.Range("R97").Value = ""
$ws.Range("S97").ClearContents()
$ws.Range("T97").ClearContents()
$ws.Range("V97").ClearContents()
$ws.Range("W97").ClearContents()
$ws.Range("X97").ClearContents()
$ws.Range("Y97").ClearContents()

# Row 98 -> becomes "A 18142-2025"
$ws.Range("A98").Value = "A 18142-2025"
$ws.Range("B98").Value = 45761.60667824074
$ws.Range("G98").Value = 1.6
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 0
$ws.Range("O98").Value = 0
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 0
$ws.Range("R98").Value = ""
$ws.Range("S98").ClearContents()
$ws.Range("T98").ClearContents()
$ws.Range("V98").ClearContents()
$ws.Range("W98").ClearContents()
$ws.Range("X98").ClearContents()
$ws.Range("Y98").ClearContents()

# Row 99 -> becomes "A 18143-2025"
$ws.Range("A99").Value = "A 18143-2025"
$ws.Range("B99").Value = 45761.60916666667
$ws.Range("G99").Value = 1.3
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("O99").Value = 0
$ws.Range("P99").Value = 0
$ws.Range("Q99").Value = 0
$ws.Range("R99").Value = ""
$ws.Range("S99").ClearContents()
$ws.Range("T99").ClearContents()
$ws.Range("V99").ClearContents()
$ws.Range("W99").ClearContents()
$ws.Range("X99").ClearContents()
$ws.Range("Y99").ClearContents()

# Row 100 -> becomes "A 54740-2022"
$ws.Range("A100").Value = "A 54740-2022"
$ws.Range("B100").Value = 44883
$ws.Range("G100").Value = 5.1
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("O100").Value = 0
$ws.Range("P100").Value = 0
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = ""
$ws.Range("S100").ClearContents()
$ws.Range("T100").ClearContents()
$ws.Range("V100").ClearContents()
$ws.Range("W100").ClearContents()
$ws.Range("X100").ClearContents()
$ws.Range("Y100").ClearContents()

# Row 101 -> becomes "A 2269-2025"
$ws.Range("A101").Value = "A 2269-2025"
$ws.Range("B101").Value = 45673.54107638889
$ws.Range("G101").Value = 2.8
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = 0
$ws.Range("Q101").Value = 0
$ws.Range("R101").Value = ""
$ws.Range("S101").ClearContents()
$ws.Range("T101").ClearContents()
$ws.Range("V101").ClearContents()
$ws.Range("W101").ClearContents()
$ws.Range("X101").ClearContents()
$ws.Range("Y101").ClearContents()

# Row 102 -> becomes "A 6500-2025"
$ws.Range("A102").Value = "A 6500-2025"
$ws.Range("B102").Value = 45699.64050925926
$ws.Range("G102").Value = 2.1
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 0
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = ""
$ws.Range("S102").ClearContents()
$ws.Range("T102").ClearContents()
$ws.Range("V102").ClearContents()
$ws.Range("W102").ClearContents()
$ws.Range("X102").ClearContents()
$ws.Range("Y102").ClearContents()

# Row 103 -> becomes "A 13250-2024"
$ws.Range("A103").Value = "A 13250-2024"
$ws.Range("B103").Value = 45386.5759375
$ws.Range("G103").Value = 0.6
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("O103").Value = 0
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 0
$ws.Range("R103").Value = ""
$ws.Range("S103").ClearContents()
$ws.Range("T103").ClearContents()
$ws.Range("V103").ClearContents()
$ws.Range("W103").ClearContents()
$ws.Range("X103").ClearContents()
$ws.Range("Y103").ClearContents()

# Row 104 -> becomes "A 3408-2024"
$ws.Range("A104").Value = "A 3408-2024"
$ws.Range("B104").Value = 45319
$ws.Range("G104").Value = 0.4
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 0
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 0
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = ""
$ws.Range("S104").ClearContents()
$ws.Range("T104").ClearContents()
$ws.Range("V104").ClearContents()
$ws.Range("W104").ClearContents()
$ws.Range("X104").ClearContents()
$ws.Range("Y104").ClearContents()

# Row 105 -> becomes "A 6491-2025"
$ws.Range("A105").Value = "A 6491-2025"
$ws.Range("B105").Value = 45699.62943287037
$ws.Range("G105").Value = 2.3
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("O105").Value = 0
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = ""
$ws.Range("S105").ClearContents()
$ws.Range("T105").ClearContents()
$ws.Range("V105").ClearContents()
$ws.Range("W105").ClearContents()
$ws.Range("X105").ClearContents()
$ws.Range("Y105").ClearContents()
